$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D8").Value = 1044600
$ws.Range("E8").Value = 908500
$ws.Range("F8").Value = 953000
$ws.Range("G8").Value = 951400
$ws.Range("H8").Value = 863600
$ws.Range("I8").Value = 714200
$ws.Range("J8").Value = 714200
$ws.Range("D9").Value = 746900
$ws.Range("E9").Value = 644200
$ws.Range("F9").Value = 683200
$ws.Range("G9").Value = 616800
$ws.Range("H9").Value = 564600
$ws.Range("I9").Value = 502400
$ws.Range("J9").Value = 569100
$ws.Range("D10").Value = 297700
$ws.Range("E10").Value = 264300
$ws.Range("F10").Value = 269800
$ws.Range("G10").Value = 334700
$ws.Range("H10").Value = 299000
$ws.Range("I10").Value = 211800
$ws.Range("J10").Value = 145200
$ws.Range("D14").Value = 6700
$ws.Range("E14").Value = 8000
$ws.Range("F14").Value = 21400
$ws.Range("G14").Value = 7300
$ws.Range("H14").Value = 11200
$ws.Range("I14").Value = 4400
$ws.Range("J14").Value = 4500
$ws.Range("D17").Value = 858700
$ws.Range("E17").Value = 691900
$ws.Range("F17").Value = 509600
$ws.Range("G17").Value = 638800
$ws.Range("H17").Value = 691100
$ws.Range("I17").Value = 559400
$ws.Range("J17").Value = 599200
$ws.Range("D18").Value = 186000
$ws.Range("E18").Value = 216500
$ws.Range("F18").Value = 443500
$ws.Range("G18").Value = 312600
$ws.Range("H18").Value = 172600
$ws.Range("I18").Value = 154800
$ws.Range("J18").Value = 115000
$ws.Range("D20").Value = -84700
$ws.Range("E20").Value = -22300
$ws.Range("F20").Value = -138600
$ws.Range("G20").Value = -12900
$ws.Range("H20").Value = -59900
$ws.Range("I20").Value = -33900
$ws.Range("J20").Value = -70300
$ws.Range("D21").Value = 233200
$ws.Range("E21").Value = 315500
$ws.Range("F21").Value = 457400
$ws.Range("G21").Value = 414200
$ws.Range("H21").Value = 228300
$ws.Range("I21").Value = 235400
$ws.Range("J21").Value = 206600
$ws.Range("D22").Value = 106000
$ws.Range("E22").Value = 111100
$ws.Range("F22").Value = 98800
$ws.Range("G22").Value = 121100
$ws.Range("H22").Value = 115600
$ws.Range("I22").Value = 124600
$ws.Range("J22").Value = 115800
$ws.Range("D23").Value = -4700
$ws.Range("E23").Value = 83000
$ws.Range("F23").Value = 206100
$ws.Range("G23").Value = 178600
$ws.Range("H23").Value = -2900
$ws.Range("I23").Value = -3700
$ws.Range("J23").Value = -71100
$ws.Range("E24").Value = 21300
$ws.Range("F24").Value = 51200
$ws.Range("G24").Value = 38200
$ws.Range("H24").Value = 9200
$ws.Range("I24").Value = 11200
$ws.Range("J24").Value = -23100
$ws.Range("D26").Value = -4700
$ws.Range("E26").Value = 61700
$ws.Range("F26").Value = 154900
$ws.Range("G26").Value = 140500
$ws.Range("H26").Value = -12200
$ws.Range("I26").Value = -14900
$ws.Range("J26").Value = -47900
$ws.Range("D27").Value = 5800
$ws.Range("E27").Value = 27300
$ws.Range("F27").Value = 83700
$ws.Range("G27").Value = 63600
$ws.Range("H27").Value = -11500
$ws.Range("I27").Value = 5400
$ws.Range("J27").Value = 14400
$ws.Range("J29").Value = -6100
$ws.Range("D32").Value = 84700
$ws.Range("E32").Value = 22300
$ws.Range("F32").Value = 138600
$ws.Range("G32").Value = 12900
$ws.Range("H32").Value = 59900
$ws.Range("I32").Value = 33900
$ws.Range("J32").Value = 70300
$ws.Range("D33").Value = 5800
$ws.Range("E33").Value = 27300
$ws.Range("F33").Value = 83700
$ws.Range("G33").Value = 63600
$ws.Range("H33").Value = -11500
$ws.Range("I33").Value = 5400
$ws.Range("J33").Value = 8300
$ws.Range("D35").Value = 5800
$ws.Range("E35").Value = 27300
$ws.Range("F35").Value = 83700
$ws.Range("G35").Value = 63600
$ws.Range("H35").Value = -11500
$ws.Range("I35").Value = 5400
$ws.Range("J35").Value = 8300
$ws.Range("D41").Value = 30100
$ws.Range("E41").Value = 67300
$ws.Range("F41").Value = 13700
$ws.Range("G41").Value = 14000
$ws.Range("H41").Value = 15200
$ws.Range("I41").Value = 625500
$ws.Range("J41").Value = 96400
$ws.Range("D42").Value = 1989600
$ws.Range("E42").Value = 2407200
$ws.Range("F42").Value = 2485600
$ws.Range("G42").Value = 1997400
$ws.Range("H42").Value = 1681200
$ws.Range("I42").Value = 1288900
$ws.Range("J42").Value = 1388300
$ws.Range("D43").Value = 604900
$ws.Range("E43").Value = 547500
$ws.Range("F43").Value = 502800
$ws.Range("G43").Value = 448000
$ws.Range("H43").Value = 445100
$ws.Range("I43").Value = 406000
$ws.Range("J43").Value = 481300
$ws.Range("D44").Value = 195300
$ws.Range("E44").Value = 165500
$ws.Range("F44").Value = 170000
$ws.Range("G44").Value = 157400
$ws.Range("H44").Value = 148900
$ws.Range("I44").Value = 153800
$ws.Range("J44").Value = 161700
$ws.Range("D45").Value = 130400
$ws.Range("E45").Value = 127300
$ws.Range("F45").Value = 169600
$ws.Range("G45").Value = 132700
$ws.Range("H45").Value = 117700
$ws.Range("I45").Value = 108400
$ws.Range("J45").Value = 120500
$ws.Range("D46").Value = 2950300
$ws.Range("E46").Value = 3314700
$ws.Range("F46").Value = 3341600
$ws.Range("G46").Value = 2749600
$ws.Range("H46").Value = 2408100
$ws.Range("I46").Value = 2582600
$ws.Range("J46").Value = 2248200
$ws.Range("D47").Value = 2453700
$ws.Range("E47").Value = 2528500
$ws.Range("F47").Value = 2561300
$ws.Range("G47").Value = 2622400
$ws.Range("H47").Value = 2594200
$ws.Range("I47").Value = 2606500
$ws.Range("J47").Value = 2536400
$ws.Range("D48").Value = 3098000
$ws.Range("E48").Value = 3040900
$ws.Range("F48").Value = 2994900
$ws.Range("G48").Value = 2904800
$ws.Range("H48").Value = 2852200
$ws.Range("I48").Value = 2801400
$ws.Range("J48").Value = 2750000
$ws.Range("D49").Value = 4339600
$ws.Range("E49").Value = 4340300
$ws.Range("F49").Value = 4351700
$ws.Range("G49").Value = 4319600
$ws.Range("H49").Value = 4338700
$ws.Range("I49").Value = 4357600
$ws.Range("J49").Value = 4386500
$ws.Range("D52").Value = 1180300
$ws.Range("E52").Value = 1008600
$ws.Range("F52").Value = 1011400
$ws.Range("G52").Value = 1090600
$ws.Range("H52").Value = 1073600
$ws.Range("I52").Value = 1009200
$ws.Range("J52").Value = 1018300
$ws.Range("D54").Value = 14021800
$ws.Range("E54").Value = 14232900
$ws.Range("F54").Value = 14261000
$ws.Range("G54").Value = 13687000
$ws.Range("H54").Value = 13266800
$ws.Range("I54").Value = 13357300
$ws.Range("J54").Value = 12939500
$ws.Range("D57").Value = 649900
$ws.Range("E57").Value = 604000
$ws.Range("F57").Value = 624000
$ws.Range("G57").Value = 561300
$ws.Range("H57").Value = 514700
$ws.Range("I57").Value = 494000
$ws.Range("J57").Value = 521100
$ws.Range("D58").Value = 715000
$ws.Range("E58").Value = 724100
$ws.Range("F58").Value = 1074800
$ws.Range("G58").Value = 896400
$ws.Range("H58").Value = 899500
$ws.Range("I58").Value = 895400
$ws.Range("J58").Value = 744600
$ws.Range("D59").Value = 502100
$ws.Range("E59").Value = 508000
$ws.Range("F59").Value = 614300
$ws.Range("G59").Value = 456500
$ws.Range("H59").Value = 416700
$ws.Range("I59").Value = 373800
$ws.Range("J59").Value = 433900
$ws.Range("D60").Value = 1867000
$ws.Range("E60").Value = 1836100
$ws.Range("F60").Value = 2313100
$ws.Range("G60").Value = 1914200
$ws.Range("H60").Value = 1830900
$ws.Range("I60").Value = 1763300
$ws.Range("J60").Value = 1699600
$ws.Range("D61").Value = 6283300
$ws.Range("E61").Value = 6238500
$ws.Range("F61").Value = 5849800
$ws.Range("G61").Value = 5961400
$ws.Range("H61").Value = 5750200
$ws.Range("I61").Value = 5823700
$ws.Range("J61").Value = 5437600
$ws.Range("D62").Value = 1700000
$ws.Range("E62").Value = 1725900
$ws.Range("F62").Value = 1724400
$ws.Range("G62").Value = 1669000
$ws.Range("H62").Value = 1670700
$ws.Range("I62").Value = 1731800
$ws.Range("J62").Value = 1697800
$ws.Range("D66").Value = 12562600
$ws.Range("E66").Value = 12660200
$ws.Range("F66").Value = 12712800
$ws.Range("G66").Value = 12023600
$ws.Range("H66").Value = 11663100
$ws.Range("I66").Value = 11750000
$ws.Range("J66").Value = 11331300
$ws.Range("D72").Value = 830000
$ws.Range("E72").Value = 842300
$ws.Range("F72").Value = 815800
$ws.Range("G72").Value = 732100
$ws.Range("H72").Value = 668500
$ws.Range("I72").Value = 696600
$ws.Range("J72").Value = 691200
$ws.Range("D76").Value = 1459200
$ws.Range("E76").Value = 1572700
$ws.Range("F76").Value = 1548200
$ws.Range("G76").Value = 1663400
$ws.Range("H76").Value = 1603700
$ws.Range("I76").Value = 1607200
$ws.Range("J76").Value = 1608100
$ws.Range("D81").Value = 5800
$ws.Range("E81").Value = 27300
$ws.Range("F81").Value = 83700
$ws.Range("G81").Value = 63600
$ws.Range("H81").Value = -11500
$ws.Range("I81").Value = 5400
$ws.Range("J81").Value = 8300
$ws.Range("D83").Value = 131900
$ws.Range("E83").Value = 121300
$ws.Range("F83").Value = 152500
$ws.Range("G83").Value = 114400
$ws.Range("H83").Value = 115600
$ws.Range("I83").Value = 114500
$ws.Range("J83").Value = 123700
$ws.Range("D89").Value = 272700
$ws.Range("E89").Value = 469300
$ws.Range("F89").Value = 236800
$ws.Range("G89").Value = 321400
$ws.Range("H89").Value = 241200
$ws.Range("I89").Value = 248700
$ws.Range("J89").Value = 161200
$ws.Range("D91").Value = -184900
$ws.Range("E91").Value = -151900
$ws.Range("F91").Value = -194700
$ws.Range("G91").Value = -145500
$ws.Range("H91").Value = -150800
$ws.Range("I91").Value = -139100
$ws.Range("J91").Value = -152400
$ws.Range("D94").Value = -156700
$ws.Range("E94").Value = 37200
$ws.Range("F94").Value = -557300
$ws.Range("G94").Value = 184400
$ws.Range("H94").Value = -504900
$ws.Range("I94").Value = -39300
$ws.Range("J94").Value = 323200
$ws.Range("D96").Value = -138800
$ws.Range("E96").Value = -9700
$ws.Range("F96").Value = -102600
$ws.Range("H96").Value = -134800
$ws.Range("I96").Value = -41400
$ws.Range("J96").Value = -121000
$ws.Range("D100").Value = -698500
$ws.Range("E100").Value = -119700
$ws.Range("F100").Value = -18400
$ws.Range("G100").Value = -26100
$ws.Range("H100").Value = -407500
$ws.Range("I100").Value = 307000
$ws.Range("J100").Value = -142000
$ws.Range("D101").Value = 62100
$ws.Range("E101").Value = 10000
$ws.Range("F101").Value = 19400
$ws.Range("H101").Value = 2700
$ws.Range("I101").Value = 7800
$ws.Range("D102").Value = -520300
$ws.Range("E102").Value = 396700
$ws.Range("F102").Value = -319600
$ws.Range("G102").Value = 478100
$ws.Range("H102").Value = -668400
$ws.Range("I102").Value = 524200
$ws.Range("J102").Value = 344000
